$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the leveled-annuity premium columns
$ws.Range("F1").Value = "annuity_level"
$ws.Range("G1").Value = "premium_leveled"

# Match the header formatting (bold, centered, bordered) used by the other headers
$ws.Range("A1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)

# Updated values for the existing columns (life_ins, life_ins_inc, premium)
$ws.Range("C2").Value = 0.03180428173371896
$ws.Range("D2").Value = 0.3782886007398258
$ws.Range("E2").Value = 3011.853682950033

$ws.Range("C3").Value = 0.01750452826003326
$ws.Range("D3").Value = 0.1892972727053334
$ws.Range("E3").Value = 1854.297335855797

$ws.Range("C4").Value = 0.03592843138828588
$ws.Range("D4").Value = 0.4191828682618172
$ws.Range("E4").Value = 3487.267585278737

# New annuity_level / premium_leveled values
$ws.Range("F2").Value = 13.92317044246741
$ws.Range("G2").Value = 216.319529764823

$ws.Range("F3").Value = 14.00441067997329
$ws.Range("G3").Value = 132.40809472315

$ws.Range("F4").Value = 13.8905593294462
$ws.Range("G4").Value = 251.0530715553101
